$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new team-record columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header row (bold/centered/bordered style)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team's record (same values repeated for every player row, 2 through 42)
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = 66   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 47   # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
